{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (const p of paragraphs.items) {\n  p.load(\"text\");\n}\nawait context.sync();\n\n// Locate the bullet that trails off after \"...Google Chat, but\" \u2014\n// that is the paragraph the author finished and then expanded with\n// several new findings bullets.\nconst marker = \"Most students have never/almost never used Google Chat, but\";\nlet target = null;\nfor (const p of paragraphs.items) {\n  if (p.text.indexOf(marker) === 0) {\n    target = p;\n    break;\n  }\n}\nif (!target) {\n  throw new Error(\"Could not find the target paragraph to edit.\");\n}\n\n// Complete the truncated sentence.\ntarget.insertText(\"Most students have never/almost never used Google Chat, but there were still many that did before its features were restricted.\", Word.InsertLocation.replace);\nawait context.sync();\n\n// Append the new survey-summary bullet points after it, in order,\n// each inheriting the same list/paragraph formatting as `target`.\ntarget = target.insertParagraph(\"About 60% of students say they would download a personalized communication app that made communicating with students and teachers easier. 10% would download it on a condition (if other people were using it, if it was secure, etc.)\", Word.InsertLocation.after);\ntarget = target.insertParagraph(\"Most students say they\u2019d download a personalized communication app for easier communication, typically on school assignments or projects. There are a significant number of people who say they wouldn\u2019t download this app because they have no need for it, or there are already pre-existing applications that do the same thing.\", Word.InsertLocation.after);\ntarget = target.insertParagraph(\"Most people are fine with the application being censored to be school-safe, but some aren\u2019t. A lot of people don\u2019t really care and/or only want certain things censored.\", Word.InsertLocation.after);\ntarget = target.insertParagraph(\"Students say they\u2019d use this hypothetical communication app for a variety of reasons, most prominently getting help on assignments and doing group assignments.\", Word.InsertLocation.after);\ntarget = target.insertParagraph(\"There aren\u2019t many useful extra remarks and there aren't many patterns. Some people think we already have alternatives while others would want this application to be state-wide and popular.\", Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the bullet that trails off after \"...Google Chat, but\" -- that is\n# the paragraph the author finished and then expanded with several new\n# findings bullets.\n$marker = \"Most students have never/almost never used Google Chat, but\"\n$targetIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text.StartsWith($marker)) {\n        $targetIndex = $i\n        break\n    }\n}\nif ($targetIndex -eq -1) {\n    throw \"Could not find the target paragraph to edit.\"\n}\n\n# Complete the truncated sentence.\n$d.Paragraphs.Item($targetIndex).Range.Text = \"Most students have never/almost never used Google Chat, but there were still many that did before its features were restricted.\"\n\n# Append the remaining new survey-summary bullets after it, each\n# inheriting the same numbered-list paragraph formatting, in order.\n$newBullets = @(\n    \"About 60% of students say they would download a personalized communication app that made communicating with students and teachers easier. 10% would download it on a condition (if other people were using it, if it was secure, etc.)\",\n    \"Most students say they\u2019d download a personalized communication app for easier communication, typically on school assignments or projects. There are a significant number of people who say they wouldn\u2019t download this app because they have no need for it, or there are already pre-existing applications that do the same thing.\",\n    \"Most people are fine with the application being censored to be school-safe, but some aren\u2019t. A lot of people don\u2019t really care and/or only want certain things censored.\",\n    \"Students say they\u2019d use this hypothetical communication app for a variety of reasons, most prominently getting help on assignments and doing group assignments.\",\n    \"There aren\u2019t many useful extra remarks and there aren't many patterns. Some people think we already have alternatives while others would want this application to be state-wide and popular.\"\n)\n\n$insertAfterIndex = $targetIndex\nforeach ($bulletText in $newBullets) {\n    $insertRange = $d.Paragraphs.Item($insertAfterIndex).Range\n    $insertRange.Collapse(0)  # wdCollapseEnd\n    $insertRange.InsertParagraphAfter()\n    $insertAfterIndex = $insertAfterIndex + 1\n    $d.Paragraphs.Item($insertAfterIndex).Range.Text = $bulletText\n}\n"}
